# ============================================================================
# Edit script: "Casos de teste para o caso de uso - Relatório de Custo.docx"
# Applies the content changes described by the commit "Casos de testes
# arrumados." (test cases tidied up).
# ============================================================================

$d = $word.ActiveDocument
$wdReplaceAll = 2
$wdFindContinue = 1
$TAB = [string][char]9

function Replace-InRange($range, $findText, $replaceText) {
    $ok = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)
    if (-not $ok) {
        Write-Host "WARNING: could not find text:" $findText
    }
    return $ok
}

function Replace-InParagraph($index, $findText, $replaceText) {
    $p = $d.Paragraphs.Item($index)
    Replace-InRange $p.Range $findText $replaceText
}

# ---------------------------------------------------------------------------
# Title (paragraph 1): merge the two runs
#   "Casos de teste para o caso de uso " + "Relatório de Custo"
# into a single run with the same combined text.
# ---------------------------------------------------------------------------
Replace-InParagraph 1 "Casos de teste para o caso de uso Relatório de Custo" "Casos de teste para o caso de uso Relatório de Custo"

Write-Host "Step 1 (title) done"

# ---------------------------------------------------------------------------
# Caso de teste 1 - Passos (paragraphs 8-14, numId 4)
# ---------------------------------------------------------------------------
Replace-InParagraph 8 "Acessa a opção “Relatório de Custos de Produção” na tela principal do sistema" "Selecione a opção “Relatórios”."
Replace-InParagraph 9 "Apresenta o filtro de data (inicial e final)" "Verifique que uma tela foi apresentada apresentando os filtros com a data inicial e final para serem inseridos"

$find10 = "Informa as datas inicial e final " + $TAB
Replace-InParagraph 10 $find10 "Informe a data inicial e final e clique no botão “Enviar” para mostrar o resultado."

Replace-InParagraph 11 "Seleciona “Gerar relatório”" "Verifique que uma tela foi apresentada listando todos os produtos que foram produzidos entre as datas informadas."
Replace-InParagraph 12 "Apresenta o relatório produtos e os custos atuais calculados com base nas ordens de produção" "Selecione a opção “Página Principal” para voltar a tela principal do sistema."
Replace-InParagraph 13 "Visualiza e fecha a tela de relatórios" "Verifique que a tela principal do sistema foi apresentada."

Write-Host "Step 2 (caso 1 passos 1-6) done"

# ---------------------------------------------------------------------------
# Caso de teste 1 - last "Passos" item (paragraph 14): the bullet item
# "Volta para a tela principal" + tab is replaced by a plain (non-list)
# paragraph that just contains a single space, indented to match the old
# list indent.
# ---------------------------------------------------------------------------
$p14 = $d.Paragraphs.Item(14)
$r14 = $p14.Range
$r14.MoveEnd(1, -1)
$r14.Text = " "
$p14.Style = $d.Styles.Item("Normal")
$p14.Format.SpaceAfter = 0
$p14.Format.LineSpacingRule = 0
$p14.Format.LeftIndent = 18

Write-Host "Step 3 (caso 1 passo 7 replaced) done"

# ---------------------------------------------------------------------------
# "Caso de teste 2: " title (paragraph 16): merge the two runs into one.
# ---------------------------------------------------------------------------
Replace-InParagraph 16 "Caso de teste 2: datas inválidas ou sem movimentação" "Caso de teste 2: datas inválidas ou sem movimentação"

Write-Host "Step 4 (caso 2 title) done"

# ---------------------------------------------------------------------------
# Caso de teste 2 - Passos (paragraphs 22-28, numId 5 -> 6)
# ---------------------------------------------------------------------------
Replace-InParagraph 22 "Acessa a opção “Relatório de Custos de Produção” na tela principal do sistema" "Selecione a opção “Relatórios”."
Replace-InParagraph 23 "Apresenta o filtro de data (inicial e final)" "Verifique que uma tela foi apresentada apresentando os filtros com a data inicial e final para serem inseridos"

Write-Host "Step 5 (caso 2 passos 1-2) done"

$find24 = "Informa as datas inicial e final " + $TAB
$replace24 = "Informe a data inicial e final inválidas e clique no botão “Enviar”, uma mensagem de erro será apresentada."
Replace-InParagraph 24 $find24 $replace24

Write-Host "Step 6 (caso 2 passo 3) done"

# Paragraph 25: "Seleciona “Gerar relatório”" -> "Feche a mensagem de erro."
# plus a (collapsed) "_GoBack" bookmark right after the new text.
Replace-InParagraph 25 "Seleciona “Gerar relatório”" "Feche a mensagem de erro."

$p25 = $d.Paragraphs.Item(25)
$r25 = $p25.Range
$r25.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $r25)

Write-Host "Step 7 (caso 2 passo 4 + bookmark) done"

Replace-InParagraph 26 "Apresenta o relatório produtos e os custos atuais calculados com base nas ordens de produção" "Selecione a opção “Página Principal” para voltar a tela principal do sistema."
Replace-InParagraph 27 "Visualiza e fecha a tela de relatórios" "Verifique que a tela principal do sistema foi apresentada."

Write-Host "Step 8 (caso 2 passos 5-6) done"
